# Applies updated HIGH/LOW/CLOSE/LTP/VOL/9:25-CLOSE figures to Sheet1
# ("added modular amount option to calc") and moves the active selection
# from D9 to D2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @{Row=2; B=2260; C=2226.65; D=2235.45; E=2233.35; F=6; G=2252.8},
    @{Row=3; B=389.55; C=382; D=384.95; E=384.1; F=31; G=389.2},
    @{Row=4; B=1595; C=1550; D=1569.9; E=1569.8; F=25; G=1588.5},
    @{Row=5; B=7634.5; C=7446; D=7499; E=7495.6; F=15; G=7611.75},
    @{Row=6; B=218.95; C=215.65; D=216.4; E=216.1; F=42; G=217.9},
    @{Row=7; B=194.6; C=190.65; D=191.8; E=191.25; F=513; G=194.1},
    @{Row=8; B=318.8; C=313.65; D=314.75; E=314.25; F=68; G=315.7},
    @{Row=9; B=594.5; C=586.5; D=591.8; E=592.45; F=26; G=588},
    @{Row=10; B=3538.95; C=3485; D=3515; E=3509.05; F=3; G=3505.8},
    @{Row=11; B=148.2; C=143.85; D=147.8; E=147.3; F=142; G=144.1},
    @{Row=12; B=1277.7; C=1267.1; D=1275; E=1273.95; F=13; G=1274.9},
    @{Row=13; B=1490.4; C=1477.2; D=1490; E=1487.25; F=209; G=1485.4},
    @{Row=14; B=491.5; C=482.35; D=483.95; E=484.55; F=54; G=483.4},
    @{Row=15; B=950.45; C=930.45; D=948.35; E=948.1; F=157; G=939.75},
    @{Row=16; B=1490.45; C=1466; D=1486.45; E=1486.4; F=28; G=1474.45},
    @{Row=17; B=1405.9; C=1395.9; D=1404.5; E=1404.3; F=22; G=1397.45},
    @{Row=18; B=633.5; C=620.55; D=630.15; E=630.35; F=37; G=621.55},
    @{Row=19; B=465.15; C=457.35; D=460.7; E=461; F=13; G=464.6},
    @{Row=20; B=1497; C=1475; D=1479.6; E=1479.75; F=19; G=1486.8},
    @{Row=21; B=264.9; C=258.9; D=264; E=263.95; F=24; G=261.65},
    @{Row=22; B=2339.1; C=2319.55; D=2323.1; E=2323.8; F=34; G=2338.7},
    @{Row=23; B=581.2; C=574.3; D=580.35; E=579.75; F=179; G=574.7},
    @{Row=24; B=654.85; C=647.05; D=652.35; E=652.3; F=3; G=649.15},
    @{Row=25; B=969.7; C=961; D=964; E=963.7; F=5; G=967.55},
    @{Row=26; B=649.9; C=643.3; D=644.8; E=645; F=65; G=645.85},
    @{Row=27; B=252; C=249.05; D=250; E=249.8; F=70; G=251.15},
    @{Row=28; B=120.2; C=118.6; D=119.75; E=119.65; F=258; G=118.8},
    @{Row=29; B=8695; C=8582.549999999999; D=8681; E=8685.75; F=2; G=8639.049999999999}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
}

$ws.Range("D2").Select()
